$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2897.3333
$ws.Range("I40").Value = 2934.75
$ws.Range("J40").Value = 2822.5
$ws.Range("K40").Value = 2934.75
$ws.Range("L40").Value = 2822.5
$ws.Range("M40").Value = -2759.75
$ws.Range("N40").Value = -3172.5

$ws.Range("H43").Value = 4220
$ws.Range("J43").Value = 4220
$ws.Range("L43").Value = 4220
$ws.Range("N43").Value = -4358

$ws.Range("H82").Value = 9737
$ws.Range("I82").Value = 5983
$ws.Range("K82").Value = 17949
$ws.Range("M82").Value = -17543

$ws.Range("H85").Value = 9737
$ws.Range("I85").Value = 5983
$ws.Range("K85").Value = 17949
$ws.Range("M85").Value = -16545

$ws.Range("H86").Value = 42713.44
$ws.Range("I86").Value = 50194.855
$ws.Range("K86").Value = 50194.855
$ws.Range("M86").Value = -49071.855

$ws.Range("H89").Value = 42713.44
$ws.Range("I89").Value = 50194.855
$ws.Range("K89").Value = 250974.275
$ws.Range("M89").Value = -245358.275

$ws.Range("H111").Value = 1823.8
$ws.Range("I111").Value = 1265.7
$ws.Range("K111").Value = 3797.1
$ws.Range("M111").Value = -730.1000000000004

$ws.Range("H135").Value = 1644.1666
$ws.Range("I135").Value = 1096.8857
$ws.Range("K135").Value = 9871.971300000001
$ws.Range("M135").Value = -7336.971300000001

$ws.Range("H138").Value = 4395.025
$ws.Range("J138").Value = 3863.2415
$ws.Range("L138").Value = 11589.7245
$ws.Range("N138").Value = -21869.7245

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 807.6111
$ws.Range("I2").Value = 773.54285
$ws.Range("K2").Value = 773.54285
$ws.Range("M2").Value = -660.54285

$ws.Range("H14").Value = 21861.834
$ws.Range("I14").Value = 338.33334
$ws.Range("J14").Value = 43385.332
$ws.Range("K14").Value = 338.33334
$ws.Range("L14").Value = 43385.332
$ws.Range("M14").Value = -163.33334
$ws.Range("N14").Value = -43735.332

$ws.Range("H32").Value = 15516.728
$ws.Range("I32").Value = 11127.25
$ws.Range("J32").Value = 27222
$ws.Range("K32").Value = 11127.25
$ws.Range("L32").Value = 27222
$ws.Range("M32").Value = -10840.25
$ws.Range("N32").Value = -27796

$ws.Range("H46").Value = 7756.375
$ws.Range("I46").Value = 2575
$ws.Range("K46").Value = 2575
$ws.Range("M46").Value = -2256

$ws.Range("H116").Value = 807.6111
$ws.Range("I116").Value = 773.54285
$ws.Range("K116").Value = 773.54285
$ws.Range("M116").Value = 1520.45715

$ws.Range("H125").Value = 127756.2
$ws.Range("J125").Value = 127756.2
$ws.Range("L125").Value = 127756.2
$ws.Range("N125").Value = -137596.2

$ws.Range("H133").Value = 59500
$ws.Range("J133").Value = 59500
$ws.Range("L133").Value = 59500
$ws.Range("N133").Value = -64560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 807.6111
$ws.Range("I3").Value = 773.54285
$ws.Range("K3").Value = 773.54285
$ws.Range("M3").Value = -659.54285

$ws.Range("H21").Value = 26239.625
$ws.Range("J21").Value = 26239.625
$ws.Range("L21").Value = 26239.625
$ws.Range("N21").Value = -26711.625

$ws.Range("H23").Value = 1797
$ws.Range("J23").Value = 1797
$ws.Range("L23").Value = 1797
$ws.Range("N23").Value = -2363

$ws.Range("H57").Value = 95389.5
$ws.Range("J57").Value = 95389.5
$ws.Range("L57").Value = 95389.5
$ws.Range("N57").Value = -96829.5

$ws.Range("H100").Value = 31378.285
$ws.Range("J100").Value = 31378.285
$ws.Range("L100").Value = 31378.285
$ws.Range("N100").Value = -33542.285

$ws.Range("H107").Value = 1737
$ws.Range("I107").Value = 1582.6666
$ws.Range("K107").Value = 1582.6666
$ws.Range("M107").Value = 337.3334

$ws.Range("H136").Value = 95389.5
$ws.Range("J136").Value = 95389.5
$ws.Range("L136").Value = 95389.5
$ws.Range("N136").Value = -105589.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1235.2
$ws.Range("I16").Value = 1024.8334
$ws.Range("K16").Value = 1024.8334
$ws.Range("M16").Value = -737.8334

$ws.Range("H22").Value = 490.8
$ws.Range("I22").Value = 312.6154
$ws.Range("K22").Value = 312.6154
$ws.Range("M22").Value = 37.38459999999998

$ws.Range("H105").Value = 2743.647
$ws.Range("I105").Value = 818.7692
$ws.Range("K105").Value = 818.7692
$ws.Range("M105").Value = 928.2308

$ws.Range("H113").Value = 1235.2
$ws.Range("I113").Value = 1024.8334
$ws.Range("K113").Value = 1024.8334
$ws.Range("M113").Value = 1145.1666

$ws.Range("H134").Value = 33590.75
$ws.Range("I134").Value = 38180.68
$ws.Range("J134").Value = 1461.25
$ws.Range("K134").Value = 114542.04
$ws.Range("L134").Value = 4383.75
$ws.Range("M134").Value = -112007.04
$ws.Range("N134").Value = -9453.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 596.41174
$ws.Range("J122").Value = 845
$ws.Range("L122").Value = 7605
$ws.Range("N122").Value = -12505

$ws.Range("H131").Value = 3581769.8
$ws.Range("I131").Value = 1957.25
$ws.Range("K131").Value = 5871.75
$ws.Range("M131").Value = -831.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 49995.668
$ws.Range("J39").Value = 49995.668
$ws.Range("L39").Value = 49995.668
$ws.Range("N39").Value = -51059.668

$ws.Range("H96").Value = 54000
$ws.Range("J96").Value = 54000
$ws.Range("L96").Value = 54000
$ws.Range("N96").Value = -59492

$ws.Range("H98").Value = 13944.25
$ws.Range("J98").Value = 13944.25
$ws.Range("L98").Value = 13944.25
$ws.Range("N98").Value = -19934.25

$ws.Range("H99").Value = 18955.166
$ws.Range("I99").Value = 8746.4
$ws.Range("K99").Value = 8746.4
$ws.Range("M99").Value = -6500.4

$ws.Range("H105").Value = 71248.664
$ws.Range("J105").Value = 71248.664
$ws.Range("L105").Value = 71248.664
$ws.Range("N105").Value = -78236.664

$ws.Range("H133").Value = 109997
$ws.Range("J133").Value = 109997
$ws.Range("L133").Value = 109997
$ws.Range("N133").Value = -120117

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 905.25
$ws.Range("I16").Value = 894.7143
$ws.Range("K16").Value = 894.7143
$ws.Range("M16").Value = -724.7143

$ws.Range("H23").Value = 8000
$ws.Range("I23").Value = 8000
$ws.Range("K23").Value = 8000
$ws.Range("M23").Value = -7770

$ws.Range("H61").Value = 2080.0952
$ws.Range("I61").Value = 2037.9445
$ws.Range("K61").Value = 2037.9445
$ws.Range("M61").Value = -1835.9445

$ws.Range("H93").Value = 1861.1818
$ws.Range("I93").Value = 1569.1428
$ws.Range("K93").Value = 1569.1428
$ws.Range("M93").Value = -321.1428000000001

$ws.Range("H113").Value = 2080.0952
$ws.Range("I113").Value = 2037.9445
$ws.Range("K113").Value = 2037.9445
$ws.Range("M113").Value = 132.0554999999999

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H133").Value = 69442
$ws.Range("J133").Value = 69442
$ws.Range("L133").Value = 69442
$ws.Range("N133").Value = -74502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 9956.286
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 11282.333
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 11282.333
$ws.Range("M23").Value = -1771
$ws.Range("N23").Value = -11740.333

$ws.Range("H107").Value = 782.2222
$ws.Range("I107").Value = 413.4
$ws.Range("J107").Value = 1243.25
$ws.Range("K107").Value = 1240.2
$ws.Range("L107").Value = 3729.75
$ws.Range("M107").Value = 679.8000000000002
$ws.Range("N107").Value = -7569.75

$ws.Range("H113").Value = 1243.9615
$ws.Range("I113").Value = 1013.1539
$ws.Range("K113").Value = 3039.4617
$ws.Range("M113").Value = -869.4616999999998

$ws.Range("H132").Value = 94870.63
$ws.Range("I132").Value = 94870.63
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 284611.89
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -282081.89
$ws.Range("N132").ClearContents()
